$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "sample"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = "sql"
$ws.Range("E7").Value = "Sample"
$ws.Range("G7").Value = "sql"

$ws.Range("H7").Select()
